# CS133JS_Lab05_Instructions-GroupC.docx edit
#
# 1. "Web App I for Group A - Roman Numeral Converter" becomes
#    "Web App I for Group C - Roman Numeral Converter", with the tail of
#    the heading (" A - Roman Numeral Converter") split into three runs:
#       " "  |  "C"  |  " - Roman Numeral Converter"
#    (all three keep the szCs=28 / single-underline formatting).
# 2. The (hidden) "_GoBack" bookmark, which in the starting document sits
#    inside the "addTask has two parameters..." paragraph, is moved to the
#    very end of the heading paragraph above (after all of its runs).

$d = $word.ActiveDocument

# --- Locate the heading paragraph by searching for its unique text.
$match = $d.Content
$found = $match.Find.Execute("Roman Numeral Converter", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Roman Numeral Converter' heading text"
}
$matchStart = $match.Start
$matchEnd = $match.End

# Find the paragraph that contains the match (Range.Paragraphs doesn't
# expand to the enclosing paragraph in this host, so walk the document's
# paragraph collection instead).
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -le $matchStart -and $cand.Range.End -ge $matchEnd) {
        $headingPara = $cand
    }
}
if ($headingPara -eq $null) {
    throw "Could not locate heading paragraph"
}

# Start of the run that holds the leading " A" (right after "...Group").
$groupEnd = $d.Content
$found2 = $groupEnd.Find.Execute("Web App I for Group", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Web App I for Group'"
}
$runStart = $groupEnd.End

# End of the paragraph's real content, i.e. just before its own
# end-of-paragraph mark.
$paraContentEnd = $headingPara.Range.End - 1

# Replace the lone "A" with "C" in place (scoped to the run's own range so
# the match can't bleed into the neighbouring run and steal its formatting).
$runRange = $d.Range($runStart, $paraContentEnd)
$replaced = $runRange.Find.Execute("A", $true, $false, $false, $false, `
    $false, $true, 1, $false, "C", 2)
if (-not $replaced) {
    throw "Could not replace 'A' with 'C' in the heading"
}

# Split the (now single) run into three runs - " ", "C", " - Roman
# Numeral Converter" - by toggling a character formatting property on
# just the "C" character and back off again; Word materialises that as
# its own run even though the end formatting is unchanged.
$cStart = $runStart + 1
$cRange = $d.Range($cStart, $cStart + 1)
$cRange.Font.Bold = $true
$cRange.Font.Bold = $false

# --- Move the "_GoBack" bookmark to the end of this heading paragraph.
# A truly zero-length range sitting exactly at a paragraph's content end
# can't be used directly, so temporarily extend the paragraph by one
# character, drop the bookmark at the (now interior) boundary, then
# remove the temporary character again. Adding a bookmark with a name
# that already exists elsewhere relocates it, so the old "_GoBack" in the
# "addTask" paragraph is automatically removed by this call too.
$tempInsert = $d.Range($paraContentEnd, $paraContentEnd)
$tempInsert.InsertAfter("X")

$bookmarkRange = $d.Range($paraContentEnd, $paraContentEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$tempDelete = $d.Range($paraContentEnd, $paraContentEnd + 1)
$tempDelete.Delete()
